$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.125
$ws.Range("C2").Value = 0.75
$ws.Range("S2").Value = 0.125
$ws.Range("P3").Value = 0.6666666666666666
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.1052631578947368
$ws.Range("F6").Value = 0.1052631578947368
$ws.Range("J6").Value = 0.3157894736842105
$ws.Range("Q6").Value = 0.1578947368421053
$ws.Range("R6").Value = 0.05263157894736842
$ws.Range("S6").Value = 0.2631578947368421
$ws.Range("F7").Value = 0.1
$ws.Range("J7").Value = 0.1
$ws.Range("Q7").Value = 0.4
$ws.Range("S7").Value = 0.4
$ws.Range("B8").Value = 0.05660377358490566
$ws.Range("D8").Value = 0.03773584905660377
$ws.Range("F8").Value = 0.07547169811320754
$ws.Range("J8").Value = 0.07547169811320754
$ws.Range("Q8").Value = 0.1320754716981132
$ws.Range("R8").Value = 0.1320754716981132
$ws.Range("S8").Value = 0.4905660377358491
$ws.Range("B9").Value = 0.05555555555555555
$ws.Range("J9").Value = 0.1666666666666667
$ws.Range("O9").Value = 0.1111111111111111
$ws.Range("Q9").Value = 0.2222222222222222
$ws.Range("R9").Value = 0.1666666666666667
$ws.Range("S9").Value = 0.2777777777777778
$ws.Range("B10").Value = 0.02409638554216868
$ws.Range("D10").Value = 0.01204819277108434
$ws.Range("F10").Value = 0.04819277108433735
$ws.Range("J10").Value = 0.144578313253012
$ws.Range("Q10").Value = 0.2048192771084337
$ws.Range("R10").Value = 0.03614457831325301
$ws.Range("S10").Value = 0.5301204819277109
$ws.Range("G11").Value = 0.09090909090909091
$ws.Range("K11").Value = 0.1818181818181818
$ws.Range("L11").Value = 0.4545454545454545
$ws.Range("S11").Value = 0.2727272727272727
$ws.Range("G12").Value = 0.8
$ws.Range("J12").Value = 0.2
$ws.Range("H15").Value = 0.1333333333333333
$ws.Range("I15").Value = 0.06666666666666667
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("O15").Value = 0.1333333333333333
$ws.Range("S15").Value = 0.2
$ws.Range("F16").Value = 0.2
$ws.Range("H16").Value = 0.2
$ws.Range("J16").Value = 0.4
$ws.Range("S16").Value = 0.2
$ws.Range("F17").Value = 0.08823529411764706
$ws.Range("H17").Value = 0.2647058823529412
$ws.Range("I17").Value = 0.08823529411764706
$ws.Range("J17").Value = 0.1764705882352941
$ws.Range("K17").Value = 0.08823529411764706
$ws.Range("M17").Value = 0.02941176470588235
$ws.Range("O17").Value = 0.08823529411764706
$ws.Range("S17").Value = 0.1764705882352941
$ws.Range("F18").Value = 0.07142857142857142
$ws.Range("H18").Value = 0.6428571428571429
$ws.Range("I18").Value = 0.07142857142857142
$ws.Range("J18").Value = 0.1428571428571428
$ws.Range("O18").Value = 0.07142857142857142
$ws.Range("F19").Value = 0.01652892561983471
$ws.Range("H19").Value = 0.2809917355371901
$ws.Range("I19").Value = 0.1074380165289256
$ws.Range("J19").Value = 0.3388429752066116
$ws.Range("K19").Value = 0.03305785123966942
$ws.Range("M19").Value = 0.03305785123966942
$ws.Range("O19").Value = 0.03305785123966942
$ws.Range("S19").Value = 0.1570247933884298
